{"js": "// Apply the README/docx stats fix for the Renaissance / JDK21 / ZGC movie-lens\n// benchmark table. The document is a single-column table; each row holds one\n// stat value. A handful of rows get their numeric text replaced, and the last\n// three rows (which previously held a full tab-separated stats line crammed\n// into one cell) get collapsed down to just their first field.\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// 0-based row index -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"825\",\n  5: \"0.00025\",\n  6: \"0.00006\",\n  7: \"0.00002\",\n  8: \"0.00005\",\n  9: \"0.00010\",\n  10: \"0.00010\",\n  11: \"0.04944\",\n  43: \"100\",\n  44: \"0.05\",\n  45: \"3829\",\n};\n\nfor (const [idxStr, newText] of Object.entries(updates)) {\n  const idx = Number(idxStr);\n  const cell = rows.items[idx].cells.items[0];\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Apply the README/docx stats fix for the Renaissance / JDK21 / ZGC movie-lens\n# benchmark table. The document is a single-column table; each row holds one\n# stat value. A handful of rows get their numeric text replaced, and the last\n# three rows (which previously held a full tab-separated stats line crammed\n# into one cell) get collapsed down to just their first field.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1-based row index -> new cell text (table has a single column).\n$updates = [ordered]@{\n  1  = \"0M\"\n  2  = \"0M\"\n  3  = \"0M\"\n  4  = \"825\"\n  6  = \"0.00025\"\n  7  = \"0.00006\"\n  8  = \"0.00002\"\n  9  = \"0.00005\"\n  10 = \"0.00010\"\n  11 = \"0.00010\"\n  12 = \"0.04944\"\n  44 = \"100\"\n  45 = \"0.05\"\n  46 = \"3829\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n  $cell = $t.Cell($rowIndex, 1)\n  $cell.Range.Text = $updates[$rowIndex]\n}\n"}
